$wb = $excel.ActiveWorkbook

# --- Sheet2: remove the two obsolete data rows (old rows 2 & 3), shifting
#     the remaining rows (old 4, 5) up to become rows 2 & 3.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Rows("2:3").Delete() | Out-Null

# Move Sheet2's selection to F2 (and make sure it is no longer the tab shown)
$ws2.Range("F2").Select() | Out-Null

# --- KETQUA (Sheet1) becomes the active/selected sheet with B15 selected
$ws1 = $wb.Worksheets.Item("KETQUA")
$ws1.Activate() | Out-Null
$ws1.Range("B15").Select() | Out-Null

# --- Shrink the hidden ExternalData_1 defined name so it only spans the
#     remaining query-table rows (A1:F3 instead of A1:F5)
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet2!ExternalData_1") {
        $n.RefersTo = "=Sheet2!`$A`$1:`$F`$3"
    }
}
